$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Real")

# --- 8th training run (1776 images, 10 epochs) results, row 10 ---
$ws.Range("B10").Value = 1776
$ws.Range("C10").Value = 0.9
$ws.Range("E10").Value = 0.1
$ws.Range("G10").Value = 10
$ws.Range("H10").Value = 163
$ws.Range("I10").Value = 30
$ws.Range("J10").Value = 18
$ws.Range("K10").Value = 0
$ws.Range("O10").Value = "runs\detect\train25"
$ws.Range("P10").Value = "runs\detect\train252"

# Accuracy/Precision columns use a tighter 2-decimal display now instead of 4
$ws.Range("M4:N10").NumberFormat = "0.00"

# Move the active selection as left in the saved file
[void]$ws.Range("N19").Select()
